$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 313, shifting the
# existing rows 313..370 down to 315..372.
$ws.Rows.Item(313).Insert()
$ws.Rows.Item(313).Insert()

# Fill the first new row (313) with new data.
$ws.Cells.Item(313, 1).Value = 4
$ws.Cells.Item(313, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(313, 3).Value = "Los Lagos"
$ws.Cells.Item(313, 4).Value = 44504
$ws.Cells.Item(313, 5).Value = 10
$ws.Cells.Item(313, 6).Value = 100112033
$ws.Cells.Item(313, 7).Value = "Lechuga"
$ws.Cells.Item(313, 8).Value = "Conconina(o)"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 80
$ws.Cells.Item(313, 11).Value = 10000
$ws.Cells.Item(313, 12).Value = 10000
$ws.Cells.Item(313, 13).Value = 10000
$ws.Cells.Item(313, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(313, 15).Value = "Región Metropolitana"
$ws.Cells.Item(313, 16).Value = 1000
$ws.Cells.Item(313, 17).Value = 10
$ws.Cells.Item(313, 18).Value = "Hortaliza"

# Fill the second new row (314) with new data.
$ws.Cells.Item(314, 1).Value = 4
$ws.Cells.Item(314, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(314, 3).Value = "Los Lagos"
$ws.Cells.Item(314, 4).Value = 44504
$ws.Cells.Item(314, 5).Value = 10
$ws.Cells.Item(314, 6).Value = 100112033
$ws.Cells.Item(314, 7).Value = "Lechuga"
$ws.Cells.Item(314, 8).Value = "Escarola"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 250
$ws.Cells.Item(314, 11).Value = 9000
$ws.Cells.Item(314, 12).Value = 10000
$ws.Cells.Item(314, 13).Value = 9400
$ws.Cells.Item(314, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(314, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(314, 16).Value = 627
$ws.Cells.Item(314, 17).Value = 15
$ws.Cells.Item(314, 18).Value = "Hortaliza"
